$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column Q (2020) into the new column R (2021)
$ws.Range("Q2").Copy($ws.Range("R2"))
$ws.Range("Q3").Copy($ws.Range("R3"))
$ws.Range("Q4").Copy($ws.Range("R4"))
$ws.Range("Q5").Copy($ws.Range("R5"))
$ws.Range("Q6").Copy($ws.Range("R6"))

# Write the new year's data
$ws.Range("R3").Value = 2021
$ws.Range("R4").Value = 233306
$ws.Range("R5").Value = 3.5
$ws.Range("R6").Value = 30.8

# Reflect the saved selection state in the sheet view
$ws.Range("Q15").Select()
